$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback: status text changed from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears (Overview status
# columns + the per-locale "Status" column).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# Hyperlink-style blue/underline font, matching the existing "a.md"/"b.md"
# link cells already on these sheets.
$linkColor = 15570276

function Set-HandbackRow($ws, $row, $handbackFile, $handbackDate) {
    $ws.Range("C$row").Value = $newStatus

    $ws.Range("I$row").Value = "a.md"
    $ws.Range("I$row").Font.Underline = $true
    $ws.Range("I$row").Font.Color = $linkColor
    $ws.Hyperlinks.Add($ws.Range("I$row"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddbad2cfa3d31c6ea564c22a29ac5df628619204/e2e/a.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md") | Out-Null

    $ws.Range("J$row").Value = $handbackFile
    $ws.Range("K$row").Value = $handbackDate
}

# ---------------------------------------------------------------------------
# zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback
# DateTime now populated for both data rows.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn 2 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-30 20:41:53"
Set-HandbackRow $wsZhCn 3 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-30 20:41:53"
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: same shape, different handback datetime/file.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe 2 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-30 20:42:01"
Set-HandbackRow $wsDeDe 3 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-30 20:42:01"
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
